$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above the current "totals" row (row 74).
#    This pushes: old row74 (totals) -> row75, old row75 (footer) -> row76.
$ws.Rows.Item(74).Insert()

# 2) Populate the new row 74 as a data row, cloning formatting from row 73
#    (the previous data row) for the bulk of the columns.
$ws.Range("A73:Q73").Copy($ws.Range("A74:Q74"))
$ws.Rows.Item(74).RowHeight = 25.5

# 3) Set the new item's data (item #68).
$ws.Range("A74").Value = 68
$ws.Range("C74").Value = "مناديل مبلله كبيره"
$ws.Range("H74").Value = "8:0"
$ws.Range("L74").Value = "0"
$ws.Range("N74").Value = "30.00"
$ws.Range("P74").Value = "30.0000"

# 4) Column Q of the data rows uses a distinct highlighted style; rebuild it
#    explicitly on Q74 to match the other rows in that column.
$q74 = $ws.Range("Q74")
$q74.NumberFormat = "#.00"
$q74.Font.Name = "Arial"
$q74.Font.Size = 9
$q74.Font.Color = 16777215
$q74.Interior.Color = 33023
$q74.HorizontalAlignment = -4108
$q74.VerticalAlignment = -4108
$q74.WrapText = $true
$q74.ShrinkToFit = $true
$q74.Value = "1:0"

# 5) Update the totals row (now row 75) to include the new item's amount.
$ws.Rows.Item(75).RowHeight = 24.75
$ws.Range("P75").Value = 4194.8149999999996

# 6) Update the footer/export-timestamp row (now row 76) with the new time.
$ws.Range("A76").Value = "Sunday, 27 July, 2025 6:17 PM"
